$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new blank rows at row 21 (pushes existing rows 21-42 down to 29-50)
$ws.Rows("21:28").Insert()

# Copy the formatting of the row just above (row 20) onto the new rows so the
# new cells reuse the same style definitions (borders/number formats) as the
# rest of the parameter table instead of creating brand-new styles.
$ws.Range("A20:E20").Copy()
$ws.Range("A21:E28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Column A (parameter symbol) -------------------------------------------------
$ws.Range("A21").Value = "konVBR1"
$ws.Range("A23").Value = "konVBN1"
$ws.Range("A25").Value = "konPlR1"
$ws.Range("A27").Value = "konPlN1"
$ws.Range("A22").Value = "koffVBR1"
$ws.Range("A24").Value = "koffVBN1"
$ws.Range("A26").Value = "koffPlR1"
$ws.Range("A28").Value = "koffPlN1"

# --- Column B (interpretation) ----------------------------------------------------
$ws.Range("B21").Value = "Association rate of VEGF-B binding to VEGFR1"
$ws.Range("B22").Value = "Dissociation rate of VEGF-B bound to VEGFR1"
$ws.Range("B23").Value = "Association rate of VEGF-B binding to NRP1"
$ws.Range("B24").Value = "Dissociation rate of VEGF-B bound to NRP1"
$ws.Range("B25").Value = "Association rate of PlGF binding to VEGFR1"
$ws.Range("B26").Value = "Dissociation rate of PlGF bound to VEGFR1"
$ws.Range("B27").Value = "Association rate of PlGF binding to NRP1"
$ws.Range("B28").Value = "Dissociation rate of PlGF bound to NRP1"

# --- Column C (value) - PlGF values left blank, to be filled after issue #3 ------
$ws.Range("C21").Value = 158000
$ws.Range("C22").Value = 0.00009145
$ws.Range("C23").Value = 126000
$ws.Range("C24").Value = 0.000194

# --- Column D (unit) ---------------------------------------------------------------
$ws.Range("D21").Value = "1/M/s"
$ws.Range("D22").Value = "1/s"
$ws.Range("D23").Value = "1/M/s"
$ws.Range("D24").Value = "1/s"
$ws.Range("D25").Value = "1/M/s"
$ws.Range("D26").Value = "1/s"
$ws.Range("D27").Value = "1/M/s"
$ws.Range("D28").Value = "1/s"

# --- Column E (reference) -----------------------------------------------------------
$ws.Range("E21").Value = "In-house data (Shobhan)"
$ws.Range("E22").Value = "In-house data (Shobhan)"
$ws.Range("E23").Value = "In-house data (Shobhan)"
$ws.Range("E24").Value = "In-house data (Shobhan)"
$ws.Range("E25").Value = "Hoffman et al., 2013"
$ws.Range("E26").Value = "Hoffman et al., 2013"
$ws.Range("E27").Value = "Hoffman et al., 2013"
$ws.Range("E28").Value = "Hoffman et al., 2013"

# Update the saved selection to match the author's final cursor position
$null = $ws.Range("C25").Select()
